$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)
$r = $ws.Range("L84")
$r.NumberFormat = "M/D/YYYY"
$r.Value = 42005
